$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12891
$ws1.Range("F5").Value = 40
$ws1.Range("F9").Value = 13021
$ws1.Range("F10").Value = 46
$ws1.Range("F11").Value = 35
$ws1.Range("F12").Value = 5287
$ws1.Range("F13").Value = 551
$ws1.Range("F18").Value = 44
$ws1.Range("F20").Value = 687
$ws1.Range("F21").Value = 2862
$ws1.Range("F22").Value = 6214
$ws1.Range("F23").Value = 1166
$ws1.Range("F24").Value = 3640

# Sheet: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 29

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12891
$ws4.Range("F5").Value = 40
$ws4.Range("F7").Value = 29
$ws4.Range("F10").Value = 13021
$ws4.Range("F11").Value = 46
$ws4.Range("F12").Value = 35
$ws4.Range("F13").Value = 5287
$ws4.Range("F14").Value = 551
$ws4.Range("F19").Value = 44
$ws4.Range("F21").Value = 687
$ws4.Range("F22").Value = 2862
$ws4.Range("F24").Value = 6214
$ws4.Range("F25").Value = 1166
$ws4.Range("F26").Value = 3640
